# Apply the commit's changes to the document:
#  1. Remove the stray "_GoBack" bookmark that sits between the "11" and
#     "-08-2017" runs in the date line.
#  2. In the last row of the table, rename the field "ACTUEEL" to
#     "HISTORIE" and reword its description, changing "...nog actueel is
#     kijkend..." to "...nog historisch is kijkend...", re-creating the
#     "_GoBack" bookmark right after the word "historisch".

$d = $word.ActiveDocument

# --- 1. Remove the old _GoBack bookmark near the date -----------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 2. ACTUEEL -> HISTORIE -------------------------------------------
$d.Content.Find.Execute("ACTUEEL", $true, $false, $false, $false, $false, `
    $true, 1, $false, "HISTORIE", 2) | Out-Null

# --- 3. Reword the description, splitting it into three runs ----------
$full = "Geeft aan of het dienstverband nog actueel is kijkend naar de begin- en einddatum contract (True/False)"

$searchRange = $d.Content
$found = $searchRange.Find.Execute($full, $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0)

if ($found) {
    $start = $searchRange.Start

    $wordIdx = $full.IndexOf("actueel")
    $wStart = $start + $wordIdx
    $wEnd = $wStart + "actueel".Length

    # Replace "actueel" with "historisch" in place.
    $wordRange = $d.Range($wStart, $wEnd)
    $wordRange.Text = "historisch"

    # Force a run break before "historisch" (toggling and restoring a
    # character property splits the run without altering its formatting).
    $beforeRange = $d.Range($start, $wStart)
    $beforeRange.Font.Bold = 1
    $beforeRange.Font.Bold = 0

    # Insert the "_GoBack" bookmark right after "historisch" -- adding a
    # bookmark also forces a run break at that position.
    $afterWordPos = $wStart + "historisch".Length
    $bmRange = $d.Range($afterWordPos, $afterWordPos)
    $d.Bookmarks.Add("_GoBack", $bmRange)
}
